$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Trim the "Tournament" sheet: the competition-key / host-key / venue-key.*
#    rows move out to the new "Properties" sheet, so delete them here.
#    (timezone and venue.1-12 / name / location stay put.)
# ---------------------------------------------------------------------
$tournament = $wb.Worksheets.Item("Tournament")

# venue-key.1 .. venue-key.12 currently occupy rows 19-30
$tournament.Rows("19:30").Delete()
# competition-key / host-key currently occupy rows 2-3
$tournament.Rows("2:3").Delete()

# restore the selection/view state recorded for Tournament after the edit
$tournament.Range("A17:XFD28").Select()

# ---------------------------------------------------------------------
# 2. Add the new "Properties" sheet right after "Colors".
# ---------------------------------------------------------------------
$colors = $wb.Worksheets.Item("Colors")
$properties = $wb.Worksheets.Add($null, $colors)
$properties.Name = "Properties"

$rows = @(
    @("key",         "value",                "notes"),
    @("competition",  "mens-world-cup",       ""),
    @("host",         "russia",               ""),
    @("timezone",     "Europe/Moscow",        ""),
    @("color.a",      "#94d9f5",              "cyan"),
    @("color.b",      "#fee289",              "yellow"),
    @("color.c",      "#f79d8f",              "red"),
    @("color.d",      "#c4e1b5",              "green"),
    @("color.e",      "#b0d0ee",              "blue"),
    @("color.f",      "#c0e4df",              "teal"),
    @("color.g",      "#fab077",              "orange"),
    @("color.h",      "#eecbef",              "purple"),
    @("venue.01",     "ru-moscow_luzhniki",   ""),
    @("venue.02",     "ru-ekaterinburg",      ""),
    @("venue.03",     "ru-saint-petersburg",  ""),
    @("venue.04",     "ru-sochi",             ""),
    @("venue.05",     "ru-kazan",             ""),
    @("venue.06",     "ru-moscow_otkrytiye",  ""),
    @("venue.07",     "ru-saransk",           ""),
    @("venue.08",     "ru-kaliningrad",       ""),
    @("venue.09",     "ru-samara",            ""),
    @("venue.10",     "ru-rostov-on-don",     ""),
    @("venue.11",     "ru-nizhny-novgorod",   ""),
    @("venue.12",     "ru-volgograd",         "")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $properties.Cells.Item($r, 1).Value = $rows[$i][0]
    $properties.Cells.Item($r, 2).Value = $rows[$i][1]
    if ($rows[$i][2] -ne "") {
        $properties.Cells.Item($r, 3).Value = $rows[$i][2]
    }
}

$properties.Columns("A:C").AutoFit() | Out-Null

$propTable = $properties.ListObjects.Add(1, $properties.Range("A1:C24"), $null, 1)
$propTable.Name = "Properties"

$properties.Range("A4:B4").Select()

# ---------------------------------------------------------------------
# 3. Colors sheet selection restores to match the edited file.
# ---------------------------------------------------------------------
$colors.Activate()
$colors.Range("B2:B9,J2:J9").Select()

# ---------------------------------------------------------------------
# 4. Leave "Properties" as the active / selected tab, matching the saved
#    workbook state (activeTab points at the new sheet).
# ---------------------------------------------------------------------
$properties.Activate()
$properties.Range("A4:B4").Select()
